# Insert a new weekly price record for "Ajo" / "Chino" (Macroferia Regional de
# Talca) right after the existing row 127. This pushes the former rows
# 128..238 down by one (to 129..239) and grows the sheet by a single row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(128).Insert()

$ws.Cells.Item(128, 1).Value = 5
$ws.Cells.Item(128, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(128, 3).Value = "Maule"
$ws.Cells.Item(128, 4).Value = 44589
$ws.Cells.Item(128, 5).Value = 7
$ws.Cells.Item(128, 6).Value = 100112003
$ws.Cells.Item(128, 7).Value = "Ajo"
$ws.Cells.Item(128, 8).Value = "Chino"
$ws.Cells.Item(128, 9).Value = "Primera"
$ws.Cells.Item(128, 10).Value = 200
$ws.Cells.Item(128, 11).Value = 20000
$ws.Cells.Item(128, 12).Value = 20000
$ws.Cells.Item(128, 13).Value = 20000
$ws.Cells.Item(128, 14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(128, 15).Value = "China"
$ws.Cells.Item(128, 16).Value = 2000
$ws.Cells.Item(128, 17).Value = 10
$ws.Cells.Item(128, 18).Value = "Hortaliza"
